# Apply the "Add files via upload" revision to the 117_2 confirmations sheet.
#
# What changed, in plain terms:
#   - The per-branch detail rows ("New nominations", "Carryover nominations",
#     "Confirmed", "Unconfirmed", "Withdrawn", "Returned to White House") used
#     to share generic labels; they are now prefixed with their branch name
#     (e.g. "     Civilian, New nominations").
#   - The old lone "Summary" section header (row 46, no value) is removed.
#   - In its place, two new summary rows are introduced: "Total new
#     nominations" (= old "Total nominations received this Session" value)
#     and "Total carryover nominations" (= old "Total nominations carried
#     over from the First Session" value), in that order.
#   - The remaining summary rows (Total confirmed / unconfirmed / withdrawn /
#     returned to the White House) keep their existing labels & values and
#     simply shift up one row once "Summary" is gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "Summary" header row (row 46). Everything below shifts up by one.
$ws.Rows.Item(46).Delete()

# 2) Re-label the per-branch detail rows so each includes its branch name.
#    Each branch header row is followed by a variable number of detail rows
#    (not every branch has Withdrawn / Returned to White House), so the
#    first/last detail row is given explicitly per branch.
$branches = @(
    @{ Name = "Civilian";       First = 7;  Last = 12 },
    @{ Name = "Other Civilian"; First = 14; Last = 18 },
    @{ Name = "Air Force";      First = 20; Last = 24 },
    @{ Name = "Army";           First = 26; Last = 30 },
    @{ Name = "Navy";           First = 32; Last = 35 },
    @{ Name = "Marine Corps";   First = 37; Last = 41 },
    @{ Name = "Space Force";    First = 43; Last = 45 }
)

foreach ($branch in $branches) {
    for ($row = $branch.First; $row -le $branch.Last; $row++) {
        $label = $ws.Cells.Item($row, 1).Value()
        if ($label.StartsWith("     ") -and -not $label.Contains(",")) {
            $suffix = $label.Substring(5)
            $ws.Cells.Item($row, 1).Value = "     " + $branch.Name + ", " + $suffix
        }
    }
}

# 3) Replace the old "Summary" row with the two new total rows (now at 46/47,
#    after the delete above shifted the old rows 47/48 up into place).
$ws.Cells.Item(46, 1).Value = "Total new nominations"
$ws.Cells.Item(46, 2).NumberFormat = "#,##0"
$ws.Cells.Item(46, 2).Value = 19678

$ws.Cells.Item(47, 1).Value = "Total carryover nominations"
$ws.Cells.Item(47, 2).NumberFormat = "#,##0"
$ws.Cells.Item(47, 2).Value = 2793

# Rows 48-51 (Total confirmed / unconfirmed / withdrawn / returned to the
# White House) already carry the correct labels and values after the shift,
# so nothing further is required there.
